$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ee6ad9eb-fec3-457e-9f42-80c74b5eb51d", "Daniel", "1377", "fmasadlkf@gmail.com"),
    @("14b55562-9898-4b6d-8036-b64846824fae", "Enrique", "2101", "enrique@gmail.com"),
    @("3102bfaf-f1d4-46f3-8745-c401a0151a78", "fkdjasç", "12", "jfhalksdjh@gmail.com"),
    @("3bfadc1d-7458-4960-bddd-28c6367c6473", "dsafds", "21", "dfasdf@gmail.com"),
    @("4b7488d0-fe4f-46cc-b519-e058bd2be3b7", "sdadf", "21", "fasdfsdf@gmial.com"),
    @("5ad53af3-4018-4fea-b2c4-00b6937ba621", "murilo", "murilo20", "muriloluiz380@gmail.com")
)

$startRow = 5
$endRow = $startRow + $data.Count - 1
$fillRange = $ws.Range("A$startRow`:D$endRow")

# Ensure numeric-looking values (e.g. "1377", "21") are kept as text, matching
# the source workbook where every cell in the table is stored as a string.
$fillRange.NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
}

# Drop back to the default (unstyled) cell style now that the text values are
# committed, so the new rows don't pick up a stray explicit style index.
$fillRange.Style = "Normal"
